# F' Architecture deck - "Updated with current flight information"
# Applies the authored diff via PowerPoint COM-interop (TextRange / Shapes).

$p = $ppt.ActivePresentation

function Set-ParaText {
    param($TextRange, [int]$Index, [string]$NewText)
    $para = $TextRange.Paragraphs($Index, 1)
    # Force a genuine content change so the engine re-writes the run as a
    # single clean <a:r> (re-using the paragraph's existing rPr) instead of
    # doing a no-op when old/new text happen to share a prefix/suffix.
    $para.Text = "`u{2022}"
    $para.Text = $NewText
}

# ---------------------------------------------------------------------
# Slide 1 - title slide: merge subtitle runs + refresh the date line
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subTr = $subtitle.TextFrame.TextRange

Set-ParaText $subTr 1 "Jet Propulsion Laboratory,"
Set-ParaText $subTr 2 "California Institute of Technology"
Set-ParaText $subTr 3 "1/18/2018"

# ---------------------------------------------------------------------
# Slide 2 - "What is F`?"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2).TextFrame.TextRange

Set-ParaText $body2 3 "Currently baselined for JPL Sphinx Leon3 Avionics SOC"
Set-ParaText $body2 5 "Uses the concept of software components"
Set-ParaText $body2 7 "Includes framework, code generators, build tools, Command/Telemetry GUI, and unit test environment"

# ---------------------------------------------------------------------
# Slide 3 - "Where is it being used?" - full rewrite of the bullet list
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$body3 = $slide3.Shapes.Item(2).TextFrame.TextRange

$lines3 = @(
    "Development",
    "Developed under JPL technology exploration task (2013)",
    "Matured under a number of JPL projects (2014-2017)",
    "Using established JPL flight processes/analysis tools",
    "Flew on RapidScat (2014-2016)",
    "Radar experiment on ISS",
    "Very stable with no reported software bugsf",
    "Flying on Asteria (Cubesat)",
    "Asteroid detection technology demonstrator",
    "In development for:",
    "Mars Helicopter Technology Development",
    "Lunar Flashlight (Cubesat)",
    "NEAScout (Cubesat)",
    "Available on GitHub",
    "Reference example can be run on Linux, MacOS, Cygwin and most embedded ARM processors (e.g. Raspberry Pi)",
    "https://github.jpl.nasa.gov/FPRIME/fprime-sw.git"
)
$body3.Text = [string]::Join("`r", $lines3)

# Outline level: 1 = top level (no pPr override), 2 = sub-bullet (lvl="1")
$levels3 = @(1,2,2,2,1,2,2,1,2,1,2,2,2,1,2,2)
# Point size per paragraph: top-level bullets 18pt, sub-bullets 16pt
$sizes3  = @(18,16,16,16,18,16,16,18,16,18,16,16,16,18,16,16)

for ($i = 1; $i -le $lines3.Count; $i++) {
    $para = $body3.Paragraphs($i, 1)
    if ($levels3[$i - 1] -eq 2) {
        $para.IndentLevel = 2
    }
    $para.Font.Size = $sizes3[$i - 1]
}

# Re-split runs around proper nouns, mirroring the authored markup (these
# substrings keep their paragraph's size, only the run boundary changes).
$body3.Paragraphs(5, 1).Characters(9, 9).Font.Size = 18     # "RapidScat"
$body3.Paragraphs(7, 1).Characters(39, 5).Font.Size = 16    # "bugsf"
$body3.Paragraphs(8, 1).Characters(11, 7).Font.Size = 18    # "Asteria"
$body3.Paragraphs(8, 1).Characters(20, 7).Font.Size = 18    # "Cubesat"
$body3.Paragraphs(12, 1).Characters(19, 7).Font.Size = 16   # "Cubesat" (Lunar Flashlight)
$body3.Paragraphs(13, 1).Characters(1, 8).Font.Size = 16    # "NEAScout"
$body3.Paragraphs(13, 1).Characters(11, 7).Font.Size = 16   # "Cubesat" (NEAScout)

# ---------------------------------------------------------------------
# Slide 4 - "F`: A Reusable Component Architecture"
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $body4 2 "Components are not dependent on other components, so can be easily reused."

# ---------------------------------------------------------------------
# Slide 5 - "F`: A Framework for quick development"
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $body5 7 "Developer specifies common patterns in simple XML."

# ---------------------------------------------------------------------
# Slide 6 - "F`: A Framework for reuse"
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$body6 = $slide6.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $body6 16 "Python-based lightweight ground system is provided with code"

# ---------------------------------------------------------------------
# Slide 9 - "F`: A Flight-ready Framework"
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$body9 = $slide9.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $body9 1 "In 2015-2016, C&DH components were taken through flight software processes"

# ---------------------------------------------------------------------
# Slide 10 - "Both Ends of the Scale"
# ---------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)

$msp430Box = $slide10.Shapes.Item(4)
$msp430Box.Width = 2916183 / 914400 * 72   # widen textbox to fit new caption
Set-ParaText $msp430Box.TextFrame.TextRange 1 "TI MSP430 Microcontroller"

$rackBox = $slide10.Shapes.Item(6)
Set-ParaText $rackBox.TextFrame.TextRange 1 "Rack Mount PC"

# ---------------------------------------------------------------------
# Slide master / layout - footer date placeholder ("Title Slide" layout)
# ---------------------------------------------------------------------
try {
    $titleLayout = $p.SlideMaster.CustomLayouts.Item(1)
    $dateShape = $titleLayout.Shapes.Item(4)
    $dateShape.TextFrame.TextRange.Text = "1/18/2018"
} catch {
    Write-Output ("Could not update slide layout date placeholder: " + $_.Exception.Message)
}

# ---------------------------------------------------------------------
# Notes master - footer date placeholder (best effort; some hosts do not
# allow edits to the notes master shapes)
# ---------------------------------------------------------------------
try {
    $notesMaster = $p.NotesMaster
    $notesDateShape = $notesMaster.Shapes.Item(2)
    $notesDateShape.TextFrame.TextRange.Text = "1/18/2018"
} catch {
    Write-Output ("Could not update notes master date placeholder: " + $_.Exception.Message)
}
